$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 1.334522
$ws.Range("H2").Value = 4.003566
$ws.Range("I2").Value = 0.4120903366177529
$ws.Range("J2").Value = 0.4120903366177529
$ws.Range("M2").Value = 0.1994653333333334
$ws.Range("N2").Value = 0.598396
$ws.Range("O2").Value = 0.01676579960230272
$ws.Range("P2").Value = 0.01676579960230271
$ws.Range("Q2").Value = 0.2661908755706667
$ws.Range("R2").Value = 2.395717880136
$ws.Range("S2").Value = 0.006909024001778714
$ws.Range("T2").Value = 0.006909024001778712

# Row 3
$ws.Range("G3").Value = 1.334522
$ws.Range("H3").Value = 4.003566
$ws.Range("I3").Value = 0.4120903366177529
$ws.Range("J3").Value = 0.4120903366177529
$ws.Range("O3").Value = 0.03203779682023726
$ws.Range("P3").Value = 0.03203779682023726
$ws.Range("Q3").Value = 0.5086646261573333
$ws.Range("R3").Value = 4.577981635416
$ws.Range("S3").Value = 0.01320246647614275
$ws.Range("T3").Value = 0.01320246647614275

# Row 4
$ws.Range("G4").Value = 1.334522
$ws.Range("H4").Value = 4.003566
$ws.Range("I4").Value = 0.4120903366177529
$ws.Range("J4").Value = 0.4120903366177529
$ws.Range("M4").Value = 0.2888043333333333
$ws.Range("N4").Value = 0.8664129999999999
$ws.Range("O4").Value = 0.02427507324719734
$ws.Range("P4").Value = 0.02427507324719734
$ws.Range("Q4").Value = 0.3854157365286666
$ws.Range("R4").Value = 3.468741628758
$ws.Range("S4").Value = 0.01000352310585816
$ws.Range("T4").Value = 0.01000352310585816

# Row 5
$ws.Range("G5").Value = 1.334522
$ws.Range("H5").Value = 4.003566
$ws.Range("I5").Value = 0.4120903366177529
$ws.Range("J5").Value = 0.4120903366177529
$ws.Range("M5").Value = 11.02772766666667
$ws.Range("N5").Value = 33.083183
$ws.Range("O5").Value = 0.9269213303302627
$ws.Range("P5").Value = 0.9269213303302626
$ws.Range("Q5").Value = 14.71674518117533
$ws.Range("R5").Value = 132.450706630578
$ws.Range("S5").Value = 0.3819753230339733
$ws.Range("T5").Value = 0.3819753230339732

# Row 6
$ws.Range("I6").Value = 0.4618070555578372
$ws.Range("J6").Value = 0.4618070555578372
$ws.Range("M6").Value = 0.1994653333333334
$ws.Range("N6").Value = 0.598396
$ws.Range("O6").Value = 0.01676579960230272
$ws.Range("P6").Value = 0.01676579960230271
$ws.Range("Q6").Value = 0.2983055256102223
$ws.Range("R6").Value = 2.684749730492
$ws.Range("S6").Value = 0.007742564548412175
$ws.Range("T6").Value = 0.007742564548412175

# Row 7
$ws.Range("I7").Value = 0.4618070555578372
$ws.Range("J7").Value = 0.4618070555578372
$ws.Range("O7").Value = 0.03203779682023726
$ws.Range("P7").Value = 0.03203779682023726
$ws.Range("S7").Value = 0.01479528061611401
$ws.Range("T7").Value = 0.01479528061611401

# Row 8
$ws.Range("I8").Value = 0.4618070555578372
$ws.Range("J8").Value = 0.4618070555578372
$ws.Range("M8").Value = 0.2888043333333333
$ws.Range("N8").Value = 0.8664129999999999
$ws.Range("O8").Value = 0.02427507324719734
$ws.Range("P8").Value = 0.02427507324719734
$ws.Range("Q8").Value = 0.4319142931445555
$ws.Range("R8").Value = 3.887228638301
$ws.Range("S8").Value = 0.01121040009973903
$ws.Range("T8").Value = 0.01121040009973903

# Row 9
$ws.Range("I9").Value = 0.4618070555578372
$ws.Range("J9").Value = 0.4618070555578372
$ws.Range("M9").Value = 11.02772766666667
$ws.Range("N9").Value = 33.083183
$ws.Range("O9").Value = 0.9269213303302627
$ws.Range("P9").Value = 0.9269213303302626
$ws.Range("Q9").Value = 16.49224977051011
$ws.Range("R9").Value = 148.430247934591
$ws.Range("S9").Value = 0.428058810293572
$ws.Range("T9").Value = 0.428058810293572

# Row 10
$ws.Range("E10").Value = 2
$ws.Range("F10").Value = 0.6666666666666666
$ws.Range("G10").Value = 0.4083733333333333
$ws.Range("H10").Value = 1.22512
$ws.Range("I10").Value = 0.1261026078244099
$ws.Range("J10").Value = 0.1261026078244099
$ws.Range("M10").Value = 0.1994653333333334
$ws.Range("N10").Value = 0.598396
$ws.Range("O10").Value = 0.01676579960230272
$ws.Range("P10").Value = 0.01676579960230271
$ws.Range("Q10").Value = 0.08145632305777778
$ws.Range("R10").Value = 0.7331069075200001
$ws.Range("S10").Value = 0.002114211052111827
$ws.Range("T10").Value = 0.002114211052111827

# Row 11
$ws.Range("E11").Value = 2
$ws.Range("F11").Value = 0.6666666666666666
$ws.Range("G11").Value = 0.4083733333333333
$ws.Range("H11").Value = 1.22512
$ws.Range("I11").Value = 0.1261026078244099
$ws.Range("J11").Value = 0.1261026078244099
$ws.Range("O11").Value = 0.03203779682023726
$ws.Range("P11").Value = 0.03203779682023726
$ws.Range("Q11").Value = 0.1556550352355555
$ws.Range("R11").Value = 1.40089531712
$ws.Range("S11").Value = 0.004040049727980505
$ws.Range("T11").Value = 0.004040049727980506

# Row 12
$ws.Range("E12").Value = 2
$ws.Range("F12").Value = 0.6666666666666666
$ws.Range("G12").Value = 0.4083733333333333
$ws.Range("H12").Value = 1.22512
$ws.Range("I12").Value = 0.1261026078244099
$ws.Range("J12").Value = 0.1261026078244099
$ws.Range("M12").Value = 0.2888043333333333
$ws.Range("N12").Value = 0.8664129999999999
$ws.Range("O12").Value = 0.02427507324719734
$ws.Range("P12").Value = 0.02427507324719734
$ws.Range("Q12").Value = 0.1179399882844444
$ws.Range("R12").Value = 1.06145989456
$ws.Range("S12").Value = 0.00306115004160015
$ws.Range("T12").Value = 0.00306115004160015

# Row 13
$ws.Range("E13").Value = 2
$ws.Range("F13").Value = 0.6666666666666666
$ws.Range("G13").Value = 0.4083733333333333
$ws.Range("H13").Value = 1.22512
$ws.Range("I13").Value = 0.1261026078244099
$ws.Range("J13").Value = 0.1261026078244099
$ws.Range("M13").Value = 11.02772766666667
$ws.Range("N13").Value = 33.083183
$ws.Range("O13").Value = 0.9269213303302627
$ws.Range("P13").Value = 0.9269213303302626
$ws.Range("Q13").Value = 4.503429906328888
$ws.Range("R13").Value = 40.53086915695999
$ws.Range("S13").Value = 0.1168871970027174
$ws.Range("T13").Value = 0.1168871970027174
